$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: update to reflect new record (previously held by row 7)
$ws.Range("D6").Value = 44186
$ws.Range("K6").Value = "Modesto"
$ws.Range("M6").Value = 55
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 20000
$ws.Range("Q6").Value = "$/bandeja 18 kilos"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1111
$ws.Range("T6").Value = 18

# Row 7: update to reflect new record (previously held by row 6, with updated figures)
$ws.Range("D7").Value = 44544
$ws.Range("K7").Value = "Castle Brite"
$ws.Range("M7").Value = 35
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 18000
$ws.Range("Q7").Value = "$/bandeja 7 kilos"
$ws.Range("R7").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S7").Value = 1000
$ws.Range("T7").Value = 18
